# Fleet age / build year workbook update (PDF edition 10, commit 1)
#
# Changes applied:
#   - D45: corrected flights total 9,088,987 -> 90,889,871
#   - B45 ("Grand Total" row): add an "age" value of "NA" (new shared string)
#   - Column D gets its own custom width (split off from the shared C:D
#     bestFit width) to comfortably fit the longer 8-digit total
#   - Active selection left on I24 (where the edit was made)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Grand Total row (row 45) corrections ---------------------------------
# Age column has no value for the Grand Total row today; the new edition
# labels it "NA".
$ws.Range("B45").Value = "NA"

# Fix the flights total - a digit was missing (9088987 -> 90889871).
$ws.Range("D45").Value = 90889871

# --- Column widths ----------------------------------------------------------
# Columns C and D used to share one bestFit width of 8. Now that D45 holds a
# longer number, give column D its own, wider, explicit width so the value
# is fully visible (column C keeps its original bestFit width of 8).
$ws.Columns.Item(4).ColumnWidth = 10.333333333333334

# --- Selection ---------------------------------------------------------------
# Leave the active cell where the edit was made.
$ws.Range("I24").Select()
